$wb = $excel.ActiveWorkbook

# ---- Sheet: New Horizons ----
$ws = $wb.Worksheets.Item('New Horizons')

$ws.Range('A54:N54').Copy() | Out-Null
$ws.Range('A55:N58').PasteSpecial(-4122) | Out-Null

$ws.Range('G8').Copy() | Out-Null
$ws.Range('G55').PasteSpecial(-4122) | Out-Null
$ws.Range('G8').Copy() | Out-Null
$ws.Range('G57').PasteSpecial(-4122) | Out-Null

$ws.Range('A55').Value = 'GENE'
$ws.Range('B55').Value = 'LILY'
$ws.Range('C55').Value = 'KIT'
$ws.Range('D55').Value = 'OLLIE'
$ws.Range('E55').Value = 'MEG'
$ws.Range('F55').Value = 'CORDELIUS'
$ws.Range('G55').Value = 'Equipo 2'
$ws.Range('H55').Value = 'SK|Ope'
$ws.Range('I55').Value = 'SK|Joker'
$ws.Range('J55').Value = 'SK|Yoshi825'
$ws.Range('K55').Value = 'HMB|BosS'
$ws.Range('L55').Value = 'HMB|Lukii'
$ws.Range('M55').Value = 'HMB|Symantec'
$ws.Range('N55').Value = '20250724T171923.000Z'
$ws.Range('A56').Value = 'BROCK'
$ws.Range('B56').Value = 'PIPER'
$ws.Range('C56').Value = 'LILY'
$ws.Range('D56').Value = 'GRAY'
$ws.Range('E56').Value = 'BELLE'
$ws.Range('F56').Value = 'KAZE'
$ws.Range('G56').Value = 'Equipo 1'
$ws.Range('H56').Value = 'SK|Ope'
$ws.Range('I56').Value = 'SK|Yoshi825'
$ws.Range('J56').Value = 'SK|Joker'
$ws.Range('K56').Value = 'HMB|BosS'
$ws.Range('L56').Value = 'HMB|Lukii'
$ws.Range('M56').Value = 'HMB|Symantec'
$ws.Range('N56').Value = '20250724T171418.000Z'
$ws.Range('A57').Value = 'BROCK'
$ws.Range('B57').Value = 'PIPER'
$ws.Range('C57').Value = 'LILY'
$ws.Range('D57').Value = 'GRAY'
$ws.Range('E57').Value = 'BELLE'
$ws.Range('F57').Value = 'KAZE'
$ws.Range('G57').Value = 'Equipo 2'
$ws.Range('H57').Value = 'SK|Ope'
$ws.Range('I57').Value = 'SK|Yoshi825'
$ws.Range('J57').Value = 'SK|Joker'
$ws.Range('K57').Value = 'HMB|BosS'
$ws.Range('L57').Value = 'HMB|Lukii'
$ws.Range('M57').Value = 'HMB|Symantec'
$ws.Range('N57').Value = '20250724T171203.000Z'
$ws.Range('A58').Value = 'BROCK'
$ws.Range('B58').Value = 'PIPER'
$ws.Range('C58').Value = 'LILY'
$ws.Range('D58').Value = 'GRAY'
$ws.Range('E58').Value = 'BELLE'
$ws.Range('F58').Value = 'KAZE'
$ws.Range('G58').Value = 'Equipo 1'
$ws.Range('H58').Value = 'SK|Ope'
$ws.Range('I58').Value = 'SK|Yoshi825'
$ws.Range('J58').Value = 'SK|Joker'
$ws.Range('K58').Value = 'HMB|BosS'
$ws.Range('L58').Value = 'HMB|Lukii'
$ws.Range('M58').Value = 'HMB|Symantec'
$ws.Range('N58').Value = '20250724T171023.000Z'

# ---- Sheet: Hot Potato ----
$ws = $wb.Worksheets.Item('Hot Potato')

$ws.Range('A65:N65').Copy() | Out-Null
$ws.Range('A66:N72').PasteSpecial(-4122) | Out-Null

$ws.Range('G6').Copy() | Out-Null
$ws.Range('G68').PasteSpecial(-4122) | Out-Null
$ws.Range('G6').Copy() | Out-Null
$ws.Range('G69').PasteSpecial(-4122) | Out-Null
$ws.Range('G6').Copy() | Out-Null
$ws.Range('G70').PasteSpecial(-4122) | Out-Null
$ws.Range('G6').Copy() | Out-Null
$ws.Range('G71').PasteSpecial(-4122) | Out-Null
$ws.Range('G6').Copy() | Out-Null
$ws.Range('G72').PasteSpecial(-4122) | Out-Null

$ws.Range('A66').Value = 'LILY'
$ws.Range('B66').Value = 'KAZE'
$ws.Range('C66').Value = 'R-T'
$ws.Range('D66').Value = 'CARL'
$ws.Range('E66').Value = 'CROW'
$ws.Range('F66').Value = 'BULL'
$ws.Range('G66').Value = 'Equipo 2'
$ws.Range('H66').Value = 'SK|Ope'
$ws.Range('I66').Value = 'SK|Joker'
$ws.Range('J66').Value = 'SK|Yoshi825'
$ws.Range('K66').Value = 'HMB|BosS'
$ws.Range('L66').Value = 'HMB|Lukii'
$ws.Range('M66').Value = 'HMB|Symantec'
$ws.Range('N66').Value = '20250724T170441.000Z'
$ws.Range('A67').Value = 'LILY'
$ws.Range('B67').Value = 'KAZE'
$ws.Range('C67').Value = 'R-T'
$ws.Range('D67').Value = 'CARL'
$ws.Range('E67').Value = 'CROW'
$ws.Range('F67').Value = 'BULL'
$ws.Range('G67').Value = 'Equipo 2'
$ws.Range('H67').Value = 'SK|Ope'
$ws.Range('I67').Value = 'SK|Joker'
$ws.Range('J67').Value = 'SK|Yoshi825'
$ws.Range('K67').Value = 'HMB|BosS'
$ws.Range('L67').Value = 'HMB|Lukii'
$ws.Range('M67').Value = 'HMB|Symantec'
$ws.Range('N67').Value = '20250724T170217.000Z'
$ws.Range('A68').Value = 'BONNIE'
$ws.Range('B68').Value = 'BERRY'
$ws.Range('C68').Value = 'CORDELIUS'
$ws.Range('D68').Value = 'RICO'
$ws.Range('E68').Value = 'CHUCK'
$ws.Range('F68').Value = 'CROW'
$ws.Range('G68').Value = 'Equipo 1'
$ws.Range('H68').Value = 'SK|Ope'
$ws.Range('I68').Value = 'SK|Joker'
$ws.Range('J68').Value = 'SK|Yoshi825'
$ws.Range('K68').Value = 'HMB|BosS'
$ws.Range('L68').Value = 'HMB|Symantec'
$ws.Range('M68').Value = 'HMB|Lukii'
$ws.Range('N68').Value = '20250724T165610.000Z'
$ws.Range('A69').Value = 'BONNIE'
$ws.Range('B69').Value = 'BERRY'
$ws.Range('C69').Value = 'CORDELIUS'
$ws.Range('D69').Value = 'RICO'
$ws.Range('E69').Value = 'CHUCK'
$ws.Range('F69').Value = 'CROW'
$ws.Range('G69').Value = 'Equipo 1'
$ws.Range('H69').Value = 'SK|Ope'
$ws.Range('I69').Value = 'SK|Joker'
$ws.Range('J69').Value = 'SK|Yoshi825'
$ws.Range('K69').Value = 'HMB|BosS'
$ws.Range('L69').Value = 'HMB|Symantec'
$ws.Range('M69').Value = 'HMB|Lukii'
$ws.Range('N69').Value = '20250724T165350.000Z'
$ws.Range('A70').Value = 'NITA'
$ws.Range('B70').Value = 'CHUCK'
$ws.Range('C70').Value = 'AMBER'
$ws.Range('D70').Value = 'BULL'
$ws.Range('E70').Value = 'CHARLIE'
$ws.Range('F70').Value = 'KAZE'
$ws.Range('G70').Value = 'Equipo 1'
$ws.Range('H70').Value = 'FUT|GeRo'
$ws.Range('I70').Value = 'FUT|Nowy297'
$ws.Range('J70').Value = 'FUT|MeOw'
$ws.Range('K70').Value = 'TH|LeNain'
$ws.Range('L70').Value = 'TH|iKaoss'
$ws.Range('M70').Value = 'TH|Zhar'
$ws.Range('N70').Value = '20250724T171633.000Z'
$ws.Range('A71').Value = 'NITA'
$ws.Range('B71').Value = 'CHUCK'
$ws.Range('C71').Value = 'AMBER'
$ws.Range('D71').Value = 'BULL'
$ws.Range('E71').Value = 'CHARLIE'
$ws.Range('F71').Value = 'KAZE'
$ws.Range('G71').Value = 'Equipo 1'
$ws.Range('H71').Value = 'FUT|GeRo'
$ws.Range('I71').Value = 'FUT|Nowy297'
$ws.Range('J71').Value = 'FUT|MeOw'
$ws.Range('K71').Value = 'TH|LeNain'
$ws.Range('L71').Value = 'TH|iKaoss'
$ws.Range('M71').Value = 'TH|Zhar'
$ws.Range('N71').Value = '20250724T171435.000Z'
$ws.Range('A72').Value = 'CHARLIE'
$ws.Range('B72').Value = 'LILY'
$ws.Range('C72').Value = 'SHADE'
$ws.Range('D72').Value = 'JESSIE'
$ws.Range('E72').Value = 'MELODIE'
$ws.Range('F72').Value = 'KAZE'
$ws.Range('G72').Value = 'Equipo 1'
$ws.Range('H72').Value = 'FUT|GeRo'
$ws.Range('I72').Value = 'FUT|Nowy297'
$ws.Range('J72').Value = 'FUT|MeOw'
$ws.Range('K72').Value = 'TH|iKaoss'
$ws.Range('L72').Value = 'TH|LeNain'
$ws.Range('M72').Value = 'TH|Zhar'
$ws.Range('N72').Value = '20250724T172202.000Z'

# ---- Sheet: Layer Cake ----
$ws = $wb.Worksheets.Item('Layer Cake')

$ws.Range('A64:N64').Copy() | Out-Null
$ws.Range('A65:N69').PasteSpecial(-4122) | Out-Null

$ws.Range('G4').Copy() | Out-Null
$ws.Range('G65').PasteSpecial(-4122) | Out-Null
$ws.Range('G4').Copy() | Out-Null
$ws.Range('G66').PasteSpecial(-4122) | Out-Null
$ws.Range('G4').Copy() | Out-Null
$ws.Range('G68').PasteSpecial(-4122) | Out-Null
$ws.Range('G4').Copy() | Out-Null
$ws.Range('G69').PasteSpecial(-4122) | Out-Null

$ws.Range('A65').Value = 'BYRON'
$ws.Range('B65').Value = 'MEG'
$ws.Range('C65').Value = 'HANK'
$ws.Range('D65').Value = 'ASH'
$ws.Range('E65').Value = 'POCO'
$ws.Range('F65').Value = 'LOU'
$ws.Range('G65').Value = 'Equipo 2'
$ws.Range('H65').Value = 'FUT|GeRo'
$ws.Range('I65').Value = 'FUT|MeOw'
$ws.Range('J65').Value = 'FUT|Nowy297'
$ws.Range('K65').Value = 'TH|LeNain'
$ws.Range('L65').Value = 'TH|iKaoss'
$ws.Range('M65').Value = 'TH|Zhar'
$ws.Range('N65').Value = '20250724T170831.000Z'
$ws.Range('A66').Value = 'BYRON'
$ws.Range('B66').Value = 'MEG'
$ws.Range('C66').Value = 'HANK'
$ws.Range('D66').Value = 'ASH'
$ws.Range('E66').Value = 'POCO'
$ws.Range('F66').Value = 'LOU'
$ws.Range('G66').Value = 'Equipo 2'
$ws.Range('H66').Value = 'FUT|GeRo'
$ws.Range('I66').Value = 'FUT|MeOw'
$ws.Range('J66').Value = 'FUT|Nowy297'
$ws.Range('K66').Value = 'TH|LeNain'
$ws.Range('L66').Value = 'TH|iKaoss'
$ws.Range('M66').Value = 'TH|Zhar'
$ws.Range('N66').Value = '20250724T170617.000Z'
$ws.Range('A67').Value = 'BYRON'
$ws.Range('B67').Value = 'MEG'
$ws.Range('C67').Value = 'HANK'
$ws.Range('D67').Value = 'ASH'
$ws.Range('E67').Value = 'POCO'
$ws.Range('F67').Value = 'LOU'
$ws.Range('G67').Value = 'Equipo 1'
$ws.Range('H67').Value = 'FUT|GeRo'
$ws.Range('I67').Value = 'FUT|MeOw'
$ws.Range('J67').Value = 'FUT|Nowy297'
$ws.Range('K67').Value = 'TH|LeNain'
$ws.Range('L67').Value = 'TH|iKaoss'
$ws.Range('M67').Value = 'TH|Zhar'
$ws.Range('N67').Value = '20250724T170356.000Z'
$ws.Range('A68').Value = 'HANK'
$ws.Range('B68').Value = 'KIT'
$ws.Range('C68').Value = 'GENE'
$ws.Range('D68').Value = 'ASH'
$ws.Range('E68').Value = 'R-T'
$ws.Range('F68').Value = 'TICK'
$ws.Range('G68').Value = 'Equipo 2'
$ws.Range('H68').Value = 'FUT|Nowy297'
$ws.Range('I68').Value = 'FUT|MeOw'
$ws.Range('J68').Value = 'FUT|GeRo'
$ws.Range('K68').Value = 'TH|LeNain'
$ws.Range('L68').Value = 'TH|iKaoss'
$ws.Range('M68').Value = 'TH|Zhar'
$ws.Range('N68').Value = '20250724T165826.000Z'
$ws.Range('A69').Value = 'HANK'
$ws.Range('B69').Value = 'KIT'
$ws.Range('C69').Value = 'GENE'
$ws.Range('D69').Value = 'ASH'
$ws.Range('E69').Value = 'R-T'
$ws.Range('F69').Value = 'TICK'
$ws.Range('G69').Value = 'Equipo 2'
$ws.Range('H69').Value = 'FUT|Nowy297'
$ws.Range('I69').Value = 'FUT|MeOw'
$ws.Range('J69').Value = 'FUT|GeRo'
$ws.Range('K69').Value = 'TH|LeNain'
$ws.Range('L69').Value = 'TH|iKaoss'
$ws.Range('M69').Value = 'TH|Zhar'
$ws.Range('N69').Value = '20250724T165702.000Z'
